$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row.
# Find the last used row based on column A (Beteckning) and update every
# row's "Förändrad" date from 45190 (2023-09-21) to 45192 (2023-09-23).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45192
}
